$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Version value from 0.4.0 to 0.7.0
$ws.Range("B3").Value = "0.7.0"

# Remove the "Jurisdiction" / "Chile" row entirely (was row 11),
# shifting all subsequent rows up by one.
$ws.Rows.Item(11).Delete()
